$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7740287184715271
$ws.Range("B1").Value = 1.057392597198486
$ws.Range("C1").Value = 3.418669462203979
$ws.Range("D1").Value = 1.466085314750671
$ws.Range("E1").Value = 1.573581695556641
